$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.325
$ws.Range("B7").Value = 5.486
$ws.Range("A9").Value = -21.752
$ws.Range("B12").Value = 5.57
$ws.Range("B14").Value = 6.468999999999999
$ws.Range("E15").Value = 16.417
$ws.Range("A18").Value = -22.025
$ws.Range("A20").Value = -20.178
$ws.Range("B26").Value = 5.762
$ws.Range("A27").Value = -21.381
$ws.Range("B27").Value = 5.62
$ws.Range("B29").Value = 6.109999999999999
$ws.Range("E33").Value = 17.286
$ws.Range("A35").Value = -19.982
$ws.Range("E35").Value = 16.522
$ws.Range("B37").Value = 8.376000000000001
$ws.Range("B38").Value = 5.991
$ws.Range("E38").Value = 16.531
$ws.Range("E43").Value = 17.124
$ws.Range("E44").Value = 16.72
$ws.Range("E47").Value = 16.766
$ws.Range("B51").Value = 5.927999999999999
$ws.Range("E51").Value = 16.793
$ws.Range("B52").Value = 5.368
$ws.Range("B55").Value = 5.705
$ws.Range("E57").Value = 16.471
$ws.Range("E63").Value = 17.673
$ws.Range("A69").Value = -21.565
$ws.Range("B69").Value = 5.723000000000001
$ws.Range("B70").Value = 5.140000000000001
$ws.Range("E70").Value = 17.523
$ws.Range("A76").Value = -20.66
$ws.Range("A78").Value = -20.242
$ws.Range("B81").Value = 6.273
$ws.Range("A82").Value = -21.997
$ws.Range("A83").Value = -20.146
$ws.Range("B83").Value = 7.354000000000001
$ws.Range("E88").Value = 16.288
$ws.Range("A93").Value = -21.811
$ws.Range("E99").Value = 16.624
$ws.Range("B102").Value = 7.231
